# edit.ps1
# Applies the CONNECTICUT_2022.xlsx "fixing network data cleaning scripts" edit:
#  1. Rename header columns (A1..D1) to short machine-friendly names.
#  2. Title-case every municipality/state name stored in columns A and B
#     (rows 2..484) — e.g. "Amatenango de la Frontera" -> "Amatenango De La Frontera".
#  3. Correct two D-column percentage values that carry a different last-bit
#     floating point rounding (rows 173 and 335).
#  4. Remove the trailing footnote/source rows (486-490) which shrinks the
#     used range down to A1:D484.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header renames
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# ---------------------------------------------------------------------------
# 2. Title-case the Spanish place names in columns A and B (data rows 2-484).
#    Excel's InvariantCulture TextInfo.ToTitleCase(ToLower(s)) matches the
#    source transformation (Python's str.title()) exactly for every value in
#    this sheet, so we can apply it generically instead of hard-coding each
#    of the ~86 affected cells.
# ---------------------------------------------------------------------------
$culture = [System.Globalization.CultureInfo]::InvariantCulture
$ti = $culture.TextInfo

$lastDataRow = 484
for ($r = 2; $r -le $lastDataRow; $r++) {
    foreach ($c in 1, 2) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($val -ne $null) {
            $newVal = $ti.ToTitleCase($val.ToLower())
            # NOTE: this runtime's -eq/-ne string operators are case-insensitive,
            # so use the ordinal .Equals() overload to detect a *real* change.
            if (-not $val.Equals($newVal, [System.StringComparison]::Ordinal)) {
                $cell.Value = $newVal
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 3. Fix the tiny floating point rounding differences on D173 and D335.
# ---------------------------------------------------------------------------
$ws.Range("D173").Value = 0.009790209790209793
$ws.Range("D335").Value = 0.009790209790209793

# ---------------------------------------------------------------------------
# 4. Drop the trailing footnote rows (486-490): sample size / source /
#    author / date lines that no longer belong in the cleaned data file.
# ---------------------------------------------------------------------------
$ws.Rows("486:490").Delete()
